$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 80
$ws.Range("E12").Value = "80/140"
